$d = $word.ActiveDocument

# --- Change 1: fix swapped prevalence stats for disfagia/disfonia ---
# The two paragraph sentences had their statistics swapped; restore the
# correct numbers/percentages on each sentence while keeping "disfagia"
# first and "disfonia" second.

$d.Content.Find.Execute(
    "A disfagia foi diagnosticada em 11 pacientes, com taxa de prevalência estimada em 23.4% (IC 95%: [12.8%, 38.4%]) no período estudado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A disfagia foi diagnosticada em 36 pacientes com taxa de prevalência estimada em 76.6% (IC 95%: [61.6%, 87.2%]) no período estudado.",
    2
)

$d.Content.Find.Execute(
    "A disfonia foi diagnosticada em 36 pacientes com taxa de prevalência estimada em 76.6% (IC 95%: [61.6%, 87.2%]) no período estudado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A disfonia foi diagnosticada em 11 pacientes, com taxa de prevalência estimada em 23.4% (IC 95%: [12.8%, 38.4%]) no período estudado.",
    2
)

# --- Change 2: simplify wording in the conclusion paragraph ---
$d.Content.Find.Execute(
    "Evidencia-se assim a necessidade de confirmação dos resultados com novos estudos, especialmente desenhados para confirmar ou refutar os resultados aqui obtidos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Evidencia-se assim a necessidade de novos estudos especialmente desenhados para confirmar ou refutar os resultados aqui obtidos.",
    2
)
